$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b.md (row 3) moves to "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 14:44:03"

# --- zh-cn detail sheet: row for b.md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "True"/"False" text would be auto-coerced to a Boolean by .Value, so copy the
# existing text "False" cell (O3) over instead, which preserves it as a text value.
$wsZhCn.Range("O3").Copy($wsZhCn.Range("F3"))
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-19 14:43:56"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0648b9a1cb95065baecf16f02f54dc7abbce9102/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cb0bcf768266bbe0d3df9584f1ecb4bcfe812d2/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de detail sheet: row for b.md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("O3").Copy($wsDeDe.Range("F3"))
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-19 14:44:03"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0648b9a1cb95065baecf16f02f54dc7abbce9102/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cb0bcf768266bbe0d3df9584f1ecb4bcfe812d2/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
